# Update "F" column (6th column) numeric values on the "展览" and "全部类型"
# worksheets to reflect newly generated output figures.

$wb = $excel.ActiveWorkbook

# Row -> (expected old value, new value) for worksheet "展览"
$sheetExhibitionChanges = @{
    3  = 2219
    4  = 94
    5  = 13311
    7  = 120
    11 = 997
    12 = 13810
    13 = 14459
    21 = 42
    23 = 114
    25 = 5507
    27 = 766
    28 = 349
    29 = 28
    30 = 105
}

# Row -> new value for worksheet "全部类型"
$sheetAllTypesChanges = @{
    3  = 2219
    4  = 94
    5  = 13311
    8  = 120
    12 = 997
    13 = 13810
    14 = 14459
    22 = 42
    24 = 114
    26 = 5507
    28 = 766
    29 = 349
    30 = 28
    31 = 105
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $sheetExhibitionChanges.Keys) {
    $wsExhibition.Range("F$row").Value = $sheetExhibitionChanges[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheetAllTypesChanges.Keys) {
    $wsAllTypes.Range("F$row").Value = $sheetAllTypesChanges[$row]
}
